$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (shared string used in A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 06:54"

# Row 6 - India
$ws.Range("B6").Value = 4113811
$ws.Range("C6").Value = 2972
$ws.Range("D6").Value = 3180865
$ws.Range("E6").Value = 862267

# Row 62 - Kirguistan
$ws.Range("B62").Value = 44316
$ws.Range("C62").Value = 23
$ws.Range("D62").Value = 39815
$ws.Range("E62").Value = 3441

# Row 125 - Tailandia
$ws.Range("B125").Value = 3444
$ws.Range("C125").Value = 6
$ws.Range("D125").Value = 3281
$ws.Range("E125").Value = 105

# Row 133 - Jamaica
$ws.Range("B133").Value = 3024
$ws.Range("C133").Value = 60
$ws.Range("D133").Value = 950
$ws.Range("E133").Value = 2042
$ws.Range("G133").Value = 2
$ws.Range("H133").Value = 32

# Row 172 - Islas Turcas y Caicos
$ws.Range("B172").Value = 598
$ws.Range("C172").Value = 7
$ws.Range("D172").Value = 264
$ws.Range("E172").Value = 329
